# CALC + CDP re-correction commit1
#
# Two students' MA2024 grades were re-corrected:
#   - Index 230261: MA2024 grade C+ -> B-  (SGPA 3.424 -> 3.458)
#   - Index 230495: MA2024 grade B  -> B+  (SGPA 3.417 -> 3.443)
#
# Those SGPA changes move both students up in the Rank-sorted results
# table (rows 88-92), which also shifts the MA2024 grade-distribution
# summary counts (column AB, rows 5-8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- MA2024 grade-distribution summary (column AB) ---------------------
$ws.Range("AB5").Value = "4(3.5%)"   # B+ row: count of B+ in MA2024
$ws.Range("AB6").Value = "1(0.9%)"   # B  row: count of B  in MA2024
$ws.Range("AB7").Value = "4(3.5%)"   # B- row: count of B- in MA2024
$ws.Range("AB8").Value = "0(0.0%)"   # C+ row: count of C+ in MA2024

# --- Re-sorted / re-ranked results rows 88-92 ---------------------------

# Row 88 -> now Index 230261 (rank 87), with corrected MA2024 grade
$ws.Range("B88").Value = "230261"
$ws.Range("D88").Value = "A-"
$ws.Range("E88").Value = "A"
$ws.Range("F88").Value = "B-"
$ws.Range("G88").Value = "A"
$ws.Range("J88").Value = "B-"
$ws.Range("K88").Value = "B-"
$ws.Range("P88").Value = 3.458

# Row 89 -> now Index 230259 (rank 88)
$ws.Range("B89").Value = "230259"
$ws.Range("C89").Value = "A"
$ws.Range("D89").Value = "B+"
$ws.Range("E89").Value = "-"
$ws.Range("F89").Value = "C+"
$ws.Range("G89").Value = "A+"
$ws.Range("I89").Value = "A-"
$ws.Range("J89").Value = "A"
$ws.Range("K89").Value = "B+"
$ws.Range("L89").Value = "-"
$ws.Range("M89").Value = "B+"
$ws.Range("P89").Value = 3.447

# Row 90 -> now Index 230495 (rank 89), with corrected MA2024 grade
$ws.Range("B90").Value = "230495"
$ws.Range("D90").Value = "A-"
$ws.Range("E90").Value = "A"
$ws.Range("F90").Value = "B-"
$ws.Range("G90").Value = "A"
$ws.Range("I90").Value = "-"
$ws.Range("J90").Value = "A-"
$ws.Range("K90").Value = "B+"
$ws.Range("L90").Value = "A+"
$ws.Range("M90").Value = "B-"
$ws.Range("P90").Value = 3.443

# Row 91 -> now Index 230473 (rank 89, tied with row above)
$ws.Range("A91").Value = 89
$ws.Range("B91").Value = "230473"
$ws.Range("C91").Value = "-"
$ws.Range("D91").Value = "A"
$ws.Range("F91").Value = "B"
$ws.Range("I91").Value = "B+"
$ws.Range("J91").Value = "A-"
$ws.Range("K91").Value = "A-"
$ws.Range("L91").Value = "B-"
$ws.Range("M91").Value = "B-"
$ws.Range("P91").Value = 3.443

# Row 92 -> now Index 230017 (rank 91, rank 90 skipped due to the tie)
$ws.Range("B92").Value = "230017"
$ws.Range("D92").Value = "B"
$ws.Range("E92").Value = "A+"
$ws.Range("F92").Value = "B"
$ws.Range("G92").Value = "A+"
$ws.Range("I92").Value = "A-"
$ws.Range("J92").Value = "B+"
$ws.Range("K92").Value = "A"
$ws.Range("L92").Value = "-"
$ws.Range("M92").Value = "B"
$ws.Range("P92").Value = 3.441
